# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (column F) counts to the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 5183
$ws1.Range("F6").Value  = 5183
$ws1.Range("F7").Value  = 132
$ws1.Range("F13").Value = 5089
$ws1.Range("F14").Value = 25
$ws1.Range("F17").Value = 237
$ws1.Range("F18").Value = 237
$ws1.Range("F21").Value = 248
$ws1.Range("F22").Value = 3842
$ws1.Range("F23").Value = 41
$ws1.Range("F24").Value = 3743
$ws1.Range("F25").Value = 181
$ws1.Range("F31").Value = 110
$ws1.Range("F37").Value = 6653
$ws1.Range("F38").Value = 1068
$ws1.Range("F39").Value = 503
$ws1.Range("F40").Value = 99
$ws1.Range("F43").Value = 1356
$ws1.Range("F45").Value = 671
$ws1.Range("F47").Value = 2272
$ws1.Range("F51").Value = 918

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 5183
$ws4.Range("F8").Value  = 5183
$ws4.Range("F9").Value  = 132
$ws4.Range("F15").Value = 5089
$ws4.Range("F16").Value = 25
$ws4.Range("F19").Value = 237
$ws4.Range("F20").Value = 237
$ws4.Range("F23").Value = 248
$ws4.Range("F24").Value = 3842
$ws4.Range("F25").Value = 3743
$ws4.Range("F26").Value = 181
$ws4.Range("F31").Value = 110
$ws4.Range("F37").Value = 6654
$ws4.Range("F38").Value = 1068
$ws4.Range("F39").Value = 503
$ws4.Range("F41").Value = 99
$ws4.Range("F44").Value = 1356
$ws4.Range("F46").Value = 671
$ws4.Range("F47").Value = 2272
$ws4.Range("F50").Value = 918
